$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Swap the match data (columns F:V) between row-pairs 11/12,
#    16/17 and 20/21. Columns A-E (index, country, tournament,
#    season, match date) are identical within each pair and stay put.
# -----------------------------------------------------------------

# swap F:V between row 11 and row 12
$ws.Cells.Item(11, 6).Value = "Al Bataeh"
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = "Ittihad Kalba"
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 3.31
$ws.Cells.Item(11, 11).Value = "21/08/2023 16:42"
$ws.Cells.Item(11, 12).Value = 3.66
$ws.Cells.Item(11, 13).Value = "26/08/2023 15:51"
$ws.Cells.Item(11, 14).Value = 3.67
$ws.Cells.Item(11, 15).Value = "21/08/2023 16:42"
$ws.Cells.Item(11, 16).Value = 3.73
$ws.Cells.Item(11, 17).Value = "26/08/2023 15:51"
$ws.Cells.Item(11, 18).Value = 2.01
$ws.Cells.Item(11, 19).Value = "21/08/2023 16:42"
$ws.Cells.Item(11, 20).Value = 1.98
$ws.Cells.Item(11, 21).Value = "26/08/2023 15:51"
$ws.Cells.Item(11, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-ittihad-kalba/678yqaxr/"
$ws.Cells.Item(12, 6).Value = "Hatta"
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = "Al Ain"
$ws.Cells.Item(12, 9).Value = 2
$ws.Cells.Item(12, 10).Value = 6.96
$ws.Cells.Item(12, 11).Value = "24/08/2023 22:39"
$ws.Cells.Item(12, 12).Value = 8.74
$ws.Cells.Item(12, 13).Value = "26/08/2023 15:57"
$ws.Cells.Item(12, 14).Value = 5.27
$ws.Cells.Item(12, 15).Value = "24/08/2023 22:39"
$ws.Cells.Item(12, 16).Value = 6.07
$ws.Cells.Item(12, 17).Value = "26/08/2023 15:57"
$ws.Cells.Item(12, 18).Value = 1.34
$ws.Cells.Item(12, 19).Value = "24/08/2023 22:39"
$ws.Cells.Item(12, 20).Value = 1.3
$ws.Cells.Item(12, 21).Value = "26/08/2023 15:57"
$ws.Cells.Item(12, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/hatta-al-ain/h0S0k1aE/"

# swap F:V between row 16 and row 17
$ws.Cells.Item(16, 6).Value = "Ittihad Kalba"
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = "Al Jazira"
$ws.Cells.Item(16, 9).Value = 4
$ws.Cells.Item(16, 10).Value = 3.33
$ws.Cells.Item(16, 11).Value = "16/09/2023 18:13"
$ws.Cells.Item(16, 12).Value = 4.23
$ws.Cells.Item(16, 13).Value = "23/09/2023 15:11"
$ws.Cells.Item(16, 14).Value = 3.77
$ws.Cells.Item(16, 15).Value = "16/09/2023 18:13"
$ws.Cells.Item(16, 16).Value = 4.42
$ws.Cells.Item(16, 17).Value = "23/09/2023 15:11"
$ws.Cells.Item(16, 18).Value = 2.04
$ws.Cells.Item(16, 19).Value = "16/09/2023 18:13"
$ws.Cells.Item(16, 20).Value = 1.71
$ws.Cells.Item(16, 21).Value = "23/09/2023 15:11"
$ws.Cells.Item(16, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/ittihad-kalba-al-jazira/x6s0ohbt/"
$ws.Cells.Item(17, 6).Value = "Al Nasr"
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = "Al Sharjah"
$ws.Cells.Item(17, 9).Value = 1
$ws.Cells.Item(17, 10).Value = 3.33
$ws.Cells.Item(17, 11).Value = "16/09/2023 18:13"
$ws.Cells.Item(17, 12).Value = 3.73
$ws.Cells.Item(17, 13).Value = "23/09/2023 15:11"
$ws.Cells.Item(17, 14).Value = 3.6
$ws.Cells.Item(17, 15).Value = "16/09/2023 18:13"
$ws.Cells.Item(17, 16).Value = 3.93
$ws.Cells.Item(17, 17).Value = "23/09/2023 15:17"
$ws.Cells.Item(17, 18).Value = 2.1
$ws.Cells.Item(17, 19).Value = "16/09/2023 18:13"
$ws.Cells.Item(17, 20).Value = 1.91
$ws.Cells.Item(17, 21).Value = "23/09/2023 15:17"
$ws.Cells.Item(17, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-nasr-al-sharjah/bytlvyiD/"

# swap F:V between row 20 and row 21
$ws.Cells.Item(20, 6).Value = "Al Bataeh"
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = "Shabab Al-Ahli Dubai"
$ws.Cells.Item(20, 9).Value = 2
$ws.Cells.Item(20, 10).Value = 6.4
$ws.Cells.Item(20, 11).Value = "17/09/2023 15:42"
$ws.Cells.Item(20, 12).Value = 6.67
$ws.Cells.Item(20, 13).Value = "24/09/2023 15:15"
$ws.Cells.Item(20, 14).Value = 5.03
$ws.Cells.Item(20, 15).Value = "17/09/2023 15:42"
$ws.Cells.Item(20, 16).Value = 4.85
$ws.Cells.Item(20, 17).Value = "24/09/2023 15:15"
$ws.Cells.Item(20, 18).Value = 1.37
$ws.Cells.Item(20, 19).Value = "17/09/2023 15:42"
$ws.Cells.Item(20, 20).Value = 1.44
$ws.Cells.Item(20, 21).Value = "24/09/2023 15:15"
$ws.Cells.Item(20, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-shabab-al-ahli-dubai/6un4pCDn/"
$ws.Cells.Item(21, 6).Value = "Al Wahda"
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = "Hatta"
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 1.36
$ws.Cells.Item(21, 11).Value = "19/09/2023 16:42"
$ws.Cells.Item(21, 12).Value = 1.24
$ws.Cells.Item(21, 13).Value = "24/09/2023 14:22"
$ws.Cells.Item(21, 14).Value = 5.48
$ws.Cells.Item(21, 15).Value = "19/09/2023 16:42"
$ws.Cells.Item(21, 16).Value = 6.54
$ws.Cells.Item(21, 17).Value = "24/09/2023 15:15"
$ws.Cells.Item(21, 18).Value = 6.99
$ws.Cells.Item(21, 19).Value = "19/09/2023 16:42"
$ws.Cells.Item(21, 20).Value = 10.02
$ws.Cells.Item(21, 21).Value = "24/09/2023 15:15"
$ws.Cells.Item(21, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wahda-hatta/pWp8qWSh/"

# -----------------------------------------------------------------
# 2) Append 7 new match rows (38..43 in football terms, rows 37..43
#    on the sheet). Copy formats from the last existing row (36)
#    first so the new rows pick up the same cell styles (bold/border
#    index column, date-formatted match-date column), then fill in
#    the values.
# -----------------------------------------------------------------
$ws.Range("A36:V36").Copy()
$ws.Range("A37:V43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 37
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "united-arab-emirates"
$ws.Cells.Item(37, 3).Value = "uae-league"
$ws.Cells.Item(37, 4).Value = "2023-2024"
$ws.Cells.Item(37, 5).Value = 45226.61458333334
$ws.Cells.Item(37, 6).Value = "Al Bataeh"
$ws.Cells.Item(37, 7).Value = 1
$ws.Cells.Item(37, 8).Value = "Ajman"
$ws.Cells.Item(37, 9).Value = 1
$ws.Cells.Item(37, 10).Value = 2.67
$ws.Cells.Item(37, 11).Value = "23/10/2023 17:42"
$ws.Cells.Item(37, 12).Value = 2.36
$ws.Cells.Item(37, 13).Value = "27/10/2023 14:40"
$ws.Cells.Item(37, 14).Value = 3.68
$ws.Cells.Item(37, 15).Value = "23/10/2023 17:42"
$ws.Cells.Item(37, 16).Value = 3.82
$ws.Cells.Item(37, 17).Value = "27/10/2023 14:36"
$ws.Cells.Item(37, 18).Value = 2.36
$ws.Cells.Item(37, 19).Value = "23/10/2023 17:42"
$ws.Cells.Item(37, 20).Value = 2.79
$ws.Cells.Item(37, 21).Value = "27/10/2023 14:40"
$ws.Cells.Item(37, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-ajman/zLLLBlCo/"

# Row 38
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "united-arab-emirates"
$ws.Cells.Item(38, 3).Value = "uae-league"
$ws.Cells.Item(38, 4).Value = "2023-2024"
$ws.Cells.Item(38, 5).Value = 45226.61458333334
$ws.Cells.Item(38, 6).Value = "Hatta"
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = "Al Wasl"
$ws.Cells.Item(38, 9).Value = 5
$ws.Cells.Item(38, 10).Value = 6.3
$ws.Cells.Item(38, 11).Value = "26/10/2023 13:42"
$ws.Cells.Item(38, 12).Value = 10.27
$ws.Cells.Item(38, 13).Value = "27/10/2023 14:13"
$ws.Cells.Item(38, 14).Value = 4.97
$ws.Cells.Item(38, 15).Value = "26/10/2023 13:42"
$ws.Cells.Item(38, 16).Value = 7.25
$ws.Cells.Item(38, 17).Value = "27/10/2023 14:13"
$ws.Cells.Item(38, 18).Value = 1.38
$ws.Cells.Item(38, 19).Value = "26/10/2023 13:42"
$ws.Cells.Item(38, 20).Value = 1.22
$ws.Cells.Item(38, 21).Value = "27/10/2023 13:29"
$ws.Cells.Item(38, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/hatta-al-wasl/SYQu7QQG/"

# Row 39
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "united-arab-emirates"
$ws.Cells.Item(39, 3).Value = "uae-league"
$ws.Cells.Item(39, 4).Value = "2023-2024"
$ws.Cells.Item(39, 5).Value = 45226.72916666666
$ws.Cells.Item(39, 6).Value = "Al Wahda"
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = "Ittihad Kalba"
$ws.Cells.Item(39, 9).Value = 1
$ws.Cells.Item(39, 10).Value = 1.62
$ws.Cells.Item(39, 11).Value = "23/10/2023 17:42"
$ws.Cells.Item(39, 12).Value = 1.71
$ws.Cells.Item(39, 13).Value = "27/10/2023 16:59"
$ws.Cells.Item(39, 14).Value = 4.52
$ws.Cells.Item(39, 15).Value = "23/10/2023 17:42"
$ws.Cells.Item(39, 16).Value = 4.5
$ws.Cells.Item(39, 17).Value = "27/10/2023 16:59"
$ws.Cells.Item(39, 18).Value = 4.16
$ws.Cells.Item(39, 19).Value = "23/10/2023 17:42"
$ws.Cells.Item(39, 20).Value = 4.11
$ws.Cells.Item(39, 21).Value = "27/10/2023 16:59"
$ws.Cells.Item(39, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wahda-ittihad-kalba/b3xwm5Zp/"

# Row 40
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "united-arab-emirates"
$ws.Cells.Item(40, 3).Value = "uae-league"
$ws.Cells.Item(40, 4).Value = "2023-2024"
$ws.Cells.Item(40, 5).Value = 45226.72916666666
$ws.Cells.Item(40, 6).Value = "Shabab Al-Ahli Dubai"
$ws.Cells.Item(40, 7).Value = 3
$ws.Cells.Item(40, 8).Value = "Al Nasr"
$ws.Cells.Item(40, 9).Value = 3
$ws.Cells.Item(40, 10).Value = 1.47
$ws.Cells.Item(40, 11).Value = "20/10/2023 17:43"
$ws.Cells.Item(40, 12).Value = 1.42
$ws.Cells.Item(40, 13).Value = "27/10/2023 17:06"
$ws.Cells.Item(40, 14).Value = 4.62
$ws.Cells.Item(40, 15).Value = "20/10/2023 17:43"
$ws.Cells.Item(40, 16).Value = 5.3
$ws.Cells.Item(40, 17).Value = "27/10/2023 17:06"
$ws.Cells.Item(40, 18).Value = 6.03
$ws.Cells.Item(40, 19).Value = "20/10/2023 17:43"
$ws.Cells.Item(40, 20).Value = 6.39
$ws.Cells.Item(40, 21).Value = "27/10/2023 17:06"
$ws.Cells.Item(40, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/shabab-al-ahli-dubai-al-nasr/8YJT9Stb/"

# Row 41
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "united-arab-emirates"
$ws.Cells.Item(41, 3).Value = "uae-league"
$ws.Cells.Item(41, 4).Value = "2023-2024"
$ws.Cells.Item(41, 5).Value = 45227.61458333334
$ws.Cells.Item(41, 6).Value = "Al Ain"
$ws.Cells.Item(41, 7).Value = 3
$ws.Cells.Item(41, 8).Value = "Emirates Club"
$ws.Cells.Item(41, 9).Value = 1
$ws.Cells.Item(41, 10).Value = 1.15
$ws.Cells.Item(41, 11).Value = "28/10/2023 12:42"
$ws.Cells.Item(41, 12).Value = 1.12
$ws.Cells.Item(41, 13).Value = "28/10/2023 14:08"
$ws.Cells.Item(41, 14).Value = 8.300000000000001
$ws.Cells.Item(41, 15).Value = "28/10/2023 12:42"
$ws.Cells.Item(41, 16).Value = 9.91
$ws.Cells.Item(41, 17).Value = "28/10/2023 14:43"
$ws.Cells.Item(41, 18).Value = 13.92
$ws.Cells.Item(41, 19).Value = "28/10/2023 12:42"
$ws.Cells.Item(41, 20).Value = 16.92
$ws.Cells.Item(41, 21).Value = "28/10/2023 14:43"
$ws.Cells.Item(41, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-ain-emirates-club/EBKPA8Ri/"

# Row 42
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "united-arab-emirates"
$ws.Cells.Item(42, 3).Value = "uae-league"
$ws.Cells.Item(42, 4).Value = "2023-2024"
$ws.Cells.Item(42, 5).Value = 45227.61458333334
$ws.Cells.Item(42, 6).Value = "Khorfakkan"
$ws.Cells.Item(42, 7).Value = 4
$ws.Cells.Item(42, 8).Value = "Al Jazira"
$ws.Cells.Item(42, 9).Value = 2
$ws.Cells.Item(42, 10).Value = 5.34
$ws.Cells.Item(42, 11).Value = "23/10/2023 17:42"
$ws.Cells.Item(42, 12).Value = 6.2
$ws.Cells.Item(42, 13).Value = "28/10/2023 14:32"
$ws.Cells.Item(42, 14).Value = 4.82
$ws.Cells.Item(42, 15).Value = "23/10/2023 17:42"
$ws.Cells.Item(42, 16).Value = 5.34
$ws.Cells.Item(42, 17).Value = "28/10/2023 14:32"
$ws.Cells.Item(42, 18).Value = 1.5
$ws.Cells.Item(42, 19).Value = "23/10/2023 17:42"
$ws.Cells.Item(42, 20).Value = 1.42
$ws.Cells.Item(42, 21).Value = "28/10/2023 14:32"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/khorfakkan-al-jazira/lrSy86BA/"

# Row 43
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "united-arab-emirates"
$ws.Cells.Item(43, 3).Value = "uae-league"
$ws.Cells.Item(43, 4).Value = "2023-2024"
$ws.Cells.Item(43, 5).Value = 45227.72916666666
$ws.Cells.Item(43, 6).Value = "Al Sharjah"
$ws.Cells.Item(43, 7).Value = 5
$ws.Cells.Item(43, 8).Value = "Bani Yas"
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 1.45
$ws.Cells.Item(43, 11).Value = "22/10/2023 21:42"
$ws.Cells.Item(43, 12).Value = 1.49
$ws.Cells.Item(43, 13).Value = "28/10/2023 17:27"
$ws.Cells.Item(43, 14).Value = 4.81
$ws.Cells.Item(43, 15).Value = "22/10/2023 21:42"
$ws.Cells.Item(43, 16).Value = 4.78
$ws.Cells.Item(43, 17).Value = "28/10/2023 17:27"
$ws.Cells.Item(43, 18).Value = 6.17
$ws.Cells.Item(43, 19).Value = "22/10/2023 21:42"
$ws.Cells.Item(43, 20).Value = 5.82
$ws.Cells.Item(43, 21).Value = "28/10/2023 17:27"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-sharjah-bani-yas/6JxX8nd4/"

Write-Host "Edit applied: swapped rows 11/12, 16/17, 20/21; appended rows 37-43."
